# Swap the "Execute" flag between the "Program Code Profile" row (D7) and
# the "PCA Group profile" row (D11), and move the active selection to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

$ws.Range("D7").Value = "no"
$ws.Range("D11").Value = "yes"

$ws.Range("D11").Select() | Out-Null
